$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.496.79"
$ws.Range("E2").Value = "  +1.26%  "
$ws.Range("D3").Value = "1.922.21"
$ws.Range("E3").Value = "  +2.01%  "
$ws.Range("D4").Value = "1.008"
$ws.Range("E4").Value = "  +0.57%  "
$ws.Range("D5").Value = "325.60"
$ws.Range("E5").Value = "  +1.07%  "
$ws.Range("D6").Value = "1.005"
$ws.Range("E6").Value = "  +0.39%  "
$ws.Range("D7").Value = "0.4838"
$ws.Range("E7").Value = "  +2.94%  "
$ws.Range("D8").Value = "0.4102"
$ws.Range("E8").Value = "  +1.95%  "
$ws.Range("D9").Value = "0.08182"
$ws.Range("E9").Value = "  +2.29%  "
$ws.Range("D10").Value = "1.025"
$ws.Range("E10").Value = "  +3.42%  "
$ws.Range("D11").Value = "23.66"
$ws.Range("E11").Value = "  +5.79%  "
$ws.Range("D12").Value = "1.958.66"
$ws.Range("E12").Value = "  +4.24%  "
$ws.Range("D13").Value = "6.064"
$ws.Range("E13").Value = "  +3.57%  "
$ws.Range("D14").Value = "7.247"
$ws.Range("E14").Value = "  +3.38%  "
$ws.Range("D15").Value = "91.50"
$ws.Range("E15").Value = "  +3.19%  "
$ws.Range("D16").Value = "0.06781"
$ws.Range("E16").Value = "  +2.72%  "
$ws.Range("E18").Value = "  +1.50%  "
$ws.Range("D19").Value = "17.83"
$ws.Range("E19").Value = "  +2.38%  "
$ws.Range("D20").Value = "1.005"
$ws.Range("E20").Value = "  +0.37%  "
$ws.Range("D21").Value = "29.534.08"
$ws.Range("E21").Value = "  +1.39%  "
$ws.Range("D22").Value = "5.636"
$ws.Range("E22").Value = "  +2.89%  "
$ws.Range("E23").Value = "  +1.42%  "
$ws.Range("D24").Value = "2.183"
$ws.Range("E24").Value = "  +0.13%  "
$ws.Range("D25").Value = "2.130.62"
$ws.Range("E25").Value = "  +0.53%  "
$ws.Range("D26").Value = "6.705"
$ws.Range("E26").Value = "  +10.66%  "
$ws.Range("D27").Value = "156.90"
$ws.Range("E27").Value = "  +1.29%  "
$ws.Range("E28").Value = "  +2.49%  "
$ws.Range("E29").Value = "  +2.69%  "
$ws.Range("D30").Value = "120.70"
$ws.Range("E30").Value = "  +2.94%  "
$ws.Range("D31").Value = "1.034"
$ws.Range("E31").Value = "  +0.32%  "
$ws.Range("E32").Value = "  +1.83%  "
$ws.Range("D33").Value = "5.540"
$ws.Range("E33").Value = "  +3.72%  "
$ws.Range("D34").Value = "3.567"
$ws.Range("E34").Value = "  +0.70%  "
$ws.Range("D35").Value = "1.393"
$ws.Range("E35").Value = "  +1.00%  "
$ws.Range("D36").Value = "0.02287"
$ws.Range("E36").Value = "  +2.63%  "
$ws.Range("D37").Value = "0.06149"
$ws.Range("E37").Value = "  +1.47%  "
$ws.Range("D38").Value = "1.180"
$ws.Range("E38").Value = "  +0.74%  "
$ws.Range("D39").Value = "0.5994"
$ws.Range("E39").Value = "  +3.43%  "
$ws.Range("D40").Value = "8.058"
$ws.Range("E40").Value = "  +0.59%  "
$ws.Range("D41").Value = "10.81"
$ws.Range("E41").Value = "  +8.11%  "
$ws.Range("E42").Value = "  +0.43%  "
$ws.Range("D43").Value = "0.1866"
$ws.Range("E43").Value = "  +2.51%  "
$ws.Range("D44").Value = "2.412"
$ws.Range("E44").Value = "  -1.98%  "
$ws.Range("D45").Value = "1.281"
$ws.Range("E45").Value = "  +0.87%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "12.54"
$ws.Range("E46").Value = "  +4.24%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "0.07609"
$ws.Range("E47").Value = "  -0.63%  "
$ws.Range("D48").Value = "0.5597"
$ws.Range("E48").Value = "  +2.58%  "
$ws.Range("D49").Value = "1.962"
$ws.Range("E49").Value = "  +3.48%  "
$ws.Range("D50").Value = "117.23"
$ws.Range("E50").Value = "  +3.33%  "
$ws.Range("E51").Value = "  +5.14%  "
